# "cambios de agosto, puntos fe de ratas e historico"
# - Update the "Nota" text in P8 (fe de erratas on the human-rights
#   recommendations boilerplate).
# - Roll the reporting period forward one quarter (B8/C8) and refresh the
#   validation/update dates (N8/O8) accordingly.
# - Minor row-height/view tidy up left behind from the author's last visit
#   to the sheet (row 3 taller header, row 8 shorter note, cursor parked
#   back on A8/A2 instead of where it had scrolled off to).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# --- Fe de erratas: replace the boilerplate "Nota" text in P8 ---------
$ws.Range("P8").Value = "La Universidad Politécnica de Pachuca, no cuenta con las recomendaciones de organismos internacionales de derechos humanos."

# --- Histórico: shift the reporting period one quarter forward --------
$ws.Range("B8").Value = 44652   # 2022-04-01 fecha de inicio
$ws.Range("C8").Value = 44742   # 2022-06-30 fecha de término
$ws.Range("N8").Value = 44753   # 2022-07-11 fecha de validación
$ws.Range("O8").Value = 44753   # 2022-07-11 fecha de actualización

# --- Row height tweaks --------------------------------------------------
$ws.Rows.Item(3).RowHeight = 23.25
$ws.Rows.Item(8).RowHeight = 58.5

# --- Restore the view / selection to the top of the table -------------
[void]$ws.Range("A8").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
